# Generate Report for Handoff
# Inserts a new localization-status entry for
#   b3843331-de7a-440b-9553-8a3392d1c054
# immediately above the existing
#   c8da94ef-7763-4898-95c5-6f57c0071fd8
# entry on every worksheet (Overview, zh-cn, de-de), pushing the
# c8da94ef row down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper data: the short xlf-content hash that is known from the diff,
# reused as a stand-in git commit SHA so the generated hyperlink target
# URLs follow the exact same scheme as every other row.
# ---------------------------------------------------------------------
$newGuid   = "b3843331-de7a-440b-9553-8a3392d1c054"
$oldGuid   = "c8da94ef-7763-4898-95c5-6f57c0071fd8"
$newCommit = "c9161aa50311a348d4412af3b7443254bdfb1ec7"

# =======================================================================
# Sheet "Overview"
# =======================================================================
$ws = $wb.Worksheets.Item("Overview")

# Duplicate row 8 (currently the c8da94ef row) down into a fresh row 8,
# pushing the original data (and its formatting) to row 9.
$ws.Rows("8:8").Copy()
$ws.Rows("8:8").Insert()

# Row 8 now becomes the new b3843331 entry.
$ws.Range("A8").Value = "$newGuid.md"
$ws.Range("B8").Value = "Ready for handoff"
$ws.Range("C8").Value = "Ready for handoff"
$ws.Range("D8").Value = "2016-36-21 00:36:08"

# Row 9 already holds the old c8da94ef values/styles from the copy, so
# nothing else needs to change there.

# Rebuild hyperlinks in final order (this engine does not shift/renumber
# hyperlink anchors on row insert, so the safest path is to clear and
# re-add every one of them top to bottom).
$ws.Range("A2").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5c1c7918917fba7188ba8b2be0407db048b3d130/e2e/a29441b6-aa6f-4c64-8ced-f0e96db3b680.md", "", "", "a29441b6-aa6f-4c64-8ced-f0e96db3b680.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/71d5fb443c3acf764e3523fad7e1e828fc1f0bdb/e2e/04f5dcfd-17e8-432b-a386-ba2204c5bc08.md", "", "", "04f5dcfd-17e8-432b-a386-ba2204c5bc08.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/7148bdbdc473d88f8fccebd584ae16b156b1c75a/e2e/1340e9cd-f921-4cbe-bb1b-d0da383a2550.md", "", "", "1340e9cd-f921-4cbe-bb1b-d0da383a2550.md")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/7148bdbdc473d88f8fccebd584ae16b156b1c75a/e2e/b899edef-7792-4edc-84d7-1435d4982e29.md", "", "", "b899edef-7792-4edc-84d7-1435d4982e29.md")
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/89dcce3a584b5b193204bd6e72e8ce3dcf027703/e2e/d9d5826d-78db-42bb-a25c-391c27bc5a40.md", "", "", "d9d5826d-78db-42bb-a25c-391c27bc5a40.md")
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/10d8aeef3c22dcce1921eec011f99d652ca29666/e2e/87295fab-7e21-42f7-81c6-2353ba2e6415.md", "", "", "87295fab-7e21-42f7-81c6-2353ba2e6415.md")
$ws.Hyperlinks.Add($ws.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/$newCommit/e2e/$newGuid.md", "", "", "$newGuid.md")
$ws.Hyperlinks.Add($ws.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/fc1ad7d9c827787ada838de4004c8d0e40c5d9f2/e2e/c8da94ef-7763-4898-95c5-6f57c0071fd8.md", "", "", "$oldGuid.md")

# =======================================================================
# Sheets "zh-cn" and "de-de" share the same layout/logic; loop over both.
# =======================================================================
$langs = @("zh-cn", "de-de")
foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang)

    $ws.Rows("8:8").Copy()
    $ws.Rows("8:8").Insert()

    if ($lang -eq "zh-cn") {
        $newXlfDate = "2016-03-21 00:36:05"
        $oldXlfDate = "2016-03-21 00:33:52"
    } else {
        $newXlfDate = "2016-03-21 00:36:08"
        $oldXlfDate = "2016-03-21 00:33:55"
    }
    $newXlf = "$newGuid.$newCommit.$lang.xlf"
    $oldXlf = "$oldGuid.196dcaf38e049fbb599ac88b9ecdbe0f263b8f2d.$lang.xlf"

    # Row 8 becomes the new b3843331 entry.
    $ws.Range("A8").Value = "$newGuid.md"
    $ws.Range("B8").Value = ".md"
    $ws.Range("C8").Value = "Ready for handoff"
    $ws.Range("D8").Value = $newXlf
    $ws.Range("E8").Value = $newXlfDate
    $ws.Range("H8").Value = "0001-01-01 00:00:00"
    $ws.Range("I8").Value = "Include"

    # Row 9 already holds the old c8da94ef row values/styles from the copy.
    # D9/E9 need the file-specific (rather than row-specific) text, which
    # the copy already preserved correctly, but set explicitly for safety.
    $ws.Range("A9").Value = "$oldGuid.md"
    $ws.Range("B9").Value = ".md"
    $ws.Range("C9").Value = "Ready for handoff"
    $ws.Range("D9").Value = $oldXlf
    $ws.Range("E9").Value = $oldXlfDate
    $ws.Range("H9").Value = "0001-01-01 00:00:00"
    $ws.Range("I9").Value = "Include"

    # Rebuild hyperlinks top to bottom in final order.
    $ws.Range("A2").Hyperlinks.Delete()

    # Row 2 (a29441b6)
    $ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5c1c7918917fba7188ba8b2be0407db048b3d130/e2e/a29441b6-aa6f-4c64-8ced-f0e96db3b680.md", "", "", "a29441b6-aa6f-4c64-8ced-f0e96db3b680.md")
    $ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/5c1c7918917fba7188ba8b2be0407db048b3d130/e2e/a29441b6-aa6f-4c64-8ced-f0e96db3b680.md", "", "", ".md")
    $ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c373222e5d6174ba69cb1a7d09d4c20894b0719b/ol-handoff/OpenLocalizationTestOrg/oltest.$lang/ci/ht/a29441b6-aa6f-4c64-8ced-f0e96db3b680.70de58809b20a8d8ab75317e2c0a9e08ebe0d72e.$lang.xlf", "", "", "a29441b6-aa6f-4c64-8ced-f0e96db3b680.70de58809b20a8d8ab75317e2c0a9e08ebe0d72e.$lang.xlf")
    $ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.$lang/blob/7a4235d8ccfb425b4b22519480915d1d90375f97/e2e/a29441b6-aa6f-4c64-8ced-f0e96db3b680.md", "", "", "a29441b6-aa6f-4c64-8ced-f0e96db3b680.md")
    $ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2b9be68e530ae9fa53f29c80c760d54eee3d786b/ol-handback/OpenLocalizationTestOrg/oltest.$lang/ci/ht/a29441b6-aa6f-4c64-8ced-f0e96db3b680.70de58809b20a8d8ab75317e2c0a9e08ebe0d72e.$lang.xlf", "", "", "a29441b6-aa6f-4c64-8ced-f0e96db3b680.70de58809b20a8d8ab75317e2c0a9e08ebe0d72e.$lang.xlf")

    # Row 3 (04f5dcfd)
    $ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/71d5fb443c3acf764e3523fad7e1e828fc1f0bdb/e2e/04f5dcfd-17e8-432b-a386-ba2204c5bc08.md", "", "", "04f5dcfd-17e8-432b-a386-ba2204c5bc08.md")
    $ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/71d5fb443c3acf764e3523fad7e1e828fc1f0bdb/e2e/04f5dcfd-17e8-432b-a386-ba2204c5bc08.md", "", "", ".md")
    $ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3c29ddf56b8f5621d5fd9061bb787a929d87c6ff/ol-handoff/OpenLocalizationTestOrg/oltest.$lang/ci/ht/04f5dcfd-17e8-432b-a386-ba2204c5bc08.75e76ab543bd9fb71f2118b17d0c30e27b7697c2.$lang.xlf", "", "", "04f5dcfd-17e8-432b-a386-ba2204c5bc08.75e76ab543bd9fb71f2118b17d0c30e27b7697c2.$lang.xlf")
    $ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.$lang/blob/f8644c7a410fe78b6fc107ed4cda9045819c8e24/e2e/04f5dcfd-17e8-432b-a386-ba2204c5bc08.md", "", "", "04f5dcfd-17e8-432b-a386-ba2204c5bc08.md")
    $ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0f2a601a4051133d010852778c72e4773c1d3ea4/ol-handback/OpenLocalizationTestOrg/oltest.$lang/ci/ht/04f5dcfd-17e8-432b-a386-ba2204c5bc08.75e76ab543bd9fb71f2118b17d0c30e27b7697c2.$lang.xlf", "", "", "04f5dcfd-17e8-432b-a386-ba2204c5bc08.75e76ab543bd9fb71f2118b17d0c30e27b7697c2.$lang.xlf")

    # Row 4 (1340e9cd)
    $ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/7148bdbdc473d88f8fccebd584ae16b156b1c75a/e2e/1340e9cd-f921-4cbe-bb1b-d0da383a2550.md", "", "", "1340e9cd-f921-4cbe-bb1b-d0da383a2550.md")
    $ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/7148bdbdc473d88f8fccebd584ae16b156b1c75a/e2e/1340e9cd-f921-4cbe-bb1b-d0da383a2550.md", "", "", ".md")
    $ws.Hyperlinks.Add($ws.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7ff54d0d2236bb5081ea43c4f0f805cff8a75001/ol-handoff/OpenLocalizationTestOrg/oltest.$lang/ci/ht/1340e9cd-f921-4cbe-bb1b-d0da383a2550.8907fc53a416d7a1b4ba9929f633c407d08a1e4f.$lang.xlf", "", "", "1340e9cd-f921-4cbe-bb1b-d0da383a2550.8907fc53a416d7a1b4ba9929f633c407d08a1e4f.$lang.xlf")

    # Row 5 (b899edef)
    $ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/7148bdbdc473d88f8fccebd584ae16b156b1c75a/e2e/b899edef-7792-4edc-84d7-1435d4982e29.md", "", "", "b899edef-7792-4edc-84d7-1435d4982e29.md")
    $ws.Hyperlinks.Add($ws.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/7148bdbdc473d88f8fccebd584ae16b156b1c75a/e2e/b899edef-7792-4edc-84d7-1435d4982e29.md", "", "", ".md")
    $ws.Hyperlinks.Add($ws.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7ff54d0d2236bb5081ea43c4f0f805cff8a75001/ol-handoff/OpenLocalizationTestOrg/oltest.$lang/ci/ht/b899edef-7792-4edc-84d7-1435d4982e29.d562657e0f8d6ab997d9a0b2844cfe0112ee6493.$lang.xlf", "", "", "b899edef-7792-4edc-84d7-1435d4982e29.d562657e0f8d6ab997d9a0b2844cfe0112ee6493.$lang.xlf")

    # Row 6 (d9d5826d)
    $ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/89dcce3a584b5b193204bd6e72e8ce3dcf027703/e2e/d9d5826d-78db-42bb-a25c-391c27bc5a40.md", "", "", "d9d5826d-78db-42bb-a25c-391c27bc5a40.md")
    $ws.Hyperlinks.Add($ws.Range("B6"), "https://github.com/OpenLocalizationTest/oltest/blob/89dcce3a584b5b193204bd6e72e8ce3dcf027703/e2e/d9d5826d-78db-42bb-a25c-391c27bc5a40.md", "", "", ".md")
    $ws.Hyperlinks.Add($ws.Range("D6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/02116fdda00689dc2b31c59332c8feac4006be7c/ol-handoff/OpenLocalizationTestOrg/oltest.$lang/ci/ht/d9d5826d-78db-42bb-a25c-391c27bc5a40.7a8ca4b25bc6ffaaf81728f3bf2a2213289309df.$lang.xlf", "", "", "d9d5826d-78db-42bb-a25c-391c27bc5a40.7a8ca4b25bc6ffaaf81728f3bf2a2213289309df.$lang.xlf")

    # Row 7 (87295fab)
    $ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/10d8aeef3c22dcce1921eec011f99d652ca29666/e2e/87295fab-7e21-42f7-81c6-2353ba2e6415.md", "", "", "87295fab-7e21-42f7-81c6-2353ba2e6415.md")
    $ws.Hyperlinks.Add($ws.Range("B7"), "https://github.com/OpenLocalizationTest/oltest/blob/10d8aeef3c22dcce1921eec011f99d652ca29666/e2e/87295fab-7e21-42f7-81c6-2353ba2e6415.md", "", "", ".md")
    $ws.Hyperlinks.Add($ws.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/93612d475839679fe25979baad56ac776bff1997/ol-handoff/OpenLocalizationTestOrg/oltest.$lang/ci/ht/87295fab-7e21-42f7-81c6-2353ba2e6415.06fcb50d84291afef3524828022ca22c3441ef2d.$lang.xlf", "", "", "87295fab-7e21-42f7-81c6-2353ba2e6415.06fcb50d84291afef3524828022ca22c3441ef2d.$lang.xlf")

    # Row 8 (NEW: b3843331)
    $ws.Hyperlinks.Add($ws.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/$newCommit/e2e/$newGuid.md", "", "", "$newGuid.md")
    $ws.Hyperlinks.Add($ws.Range("B8"), "https://github.com/OpenLocalizationTest/oltest/blob/$newCommit/e2e/$newGuid.md", "", "", ".md")
    $ws.Hyperlinks.Add($ws.Range("D8"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$newCommit/ol-handoff/OpenLocalizationTestOrg/oltest.$lang/ci/ht/$newXlf", "", "", "$newXlf")

    # Row 9 (old row 8: c8da94ef)
    $ws.Hyperlinks.Add($ws.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/fc1ad7d9c827787ada838de4004c8d0e40c5d9f2/e2e/c8da94ef-7763-4898-95c5-6f57c0071fd8.md", "", "", "$oldGuid.md")
    $ws.Hyperlinks.Add($ws.Range("B9"), "https://github.com/OpenLocalizationTest/oltest/blob/fc1ad7d9c827787ada838de4004c8d0e40c5d9f2/e2e/c8da94ef-7763-4898-95c5-6f57c0071fd8.md", "", "", ".md")
    $ws.Hyperlinks.Add($ws.Range("D9"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bff0c0cccd2aaa626e7063bb0335df354cfffbda/ol-handoff/OpenLocalizationTestOrg/oltest.$lang/ci/ht/$oldXlf", "", "", "$oldXlf")
}

$wb.Save()
